$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlink on A2 (was linking to file://signature.png)
$ws.Hyperlinks.Delete()

# Set A2 to a numeric signing code value instead of the hyperlink text
$ws.Range("A2").Value = 123123

# Apply a number format to A2 matching numFmtId 11 (0.00E+00 scientific notation)
$ws.Range("A2").NumberFormat = "0.00E+00"

# Set column A width to fit content (width 12)
$ws.Columns("A:A").ColumnWidth = 11.166666666666666
